$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424663636178707
$ws.Range("D2").Value = 0.01491041935264548
$ws.Range("E2").Value = 0.4207449291862986
$ws.Range("F2").Value = 0.5905521328082273
$ws.Range("G2").Value = 0.002395731788579087
$ws.Range("K2").Value = 0.7242949027853456
$ws.Range("N2").Value = 1.212530402582914
$ws.Range("O2").Value = 1.940657304385979
$ws.Range("B3").Value = 0.1329234035090678
$ws.Range("D3").Value = 0.01373770045224632
$ws.Range("E3").Value = 0.3670268358805089
$ws.Range("F3").Value = 0.5818197447800699
$ws.Range("G3").Value = 0.002398633927818561
$ws.Range("K3").Value = 0.6314538395737088
$ws.Range("N3").Value = 1.226253443470796
$ws.Range("O3").Value = 1.925593240060948
$ws.Range("B4").Value = 0.1271329048743723
$ws.Range("D4").Value = 0.01301234562696862
$ws.Range("E4").Value = 0.334135301463462
$ws.Range("F4").Value = 0.5769107425604005
$ws.Range("G4").Value = 0.002400509463204303
$ws.Range("K4").Value = 0.5742786145162881
$ws.Range("N4").Value = 1.235182105382741
$ws.Range("O4").Value = 1.917839695803622
$ws.Range("B5").Value = 0.1247907298525917
$ws.Range("D5").Value = 0.01271544816348325
$ws.Range("E5").Value = 0.3207529593043859
$ws.Range("F5").Value = 0.5750237985167033
$ws.Range("G5").Value = 0.002401297373906194
$ws.Range("K5").Value = 0.5509370244836589
$ws.Range("N5").Value = 1.238947024862188
$ws.Range("O5").Value = 1.915055013969408
$ws.Range("B6").Value = 0.124402875456866
$ws.Range("D6").Value = 0.01266607016869514
$ws.Range("E6").Value = 0.3185320609508153
$ws.Range("F6").Value = 0.5747173193886468
$ws.Range("G6").Value = 0.002401429634363976
$ws.Range("K6").Value = 0.5470586400951447
$ws.Range("N6").Value = 1.239579822810637
$ws.Range("O6").Value = 1.914615230967598
$ws.Range("B7").Value = 0.1271012463983681
$ws.Range("D7").Value = 0.01300834683549823
$ws.Range("E7").Value = 0.3339547389450388
$ws.Range("F7").Value = 0.5768848353923417
$ws.Range("G7").Value = 0.002400519993404245
$ws.Range("K7").Value = 0.5739639917732973
$ws.Range("N7").Value = 1.235232368520631
$ws.Range("O7").Value = 1.917800624042854
$ws.Range("B8").Value = 0.1391617460885612
$ws.Range("D8").Value = 0.01450717968332071
$ws.Range("E8").Value = 0.4022028158448592
$ws.Range("F8").Value = 0.5874470425690461
$ws.Range("G8").Value = 0.002396713063296419
$ws.Range("K8").Value = 0.6923190306126799
$ws.Range("N8").Value = 1.217157792024882
$ws.Range("O8").Value = 1.935152020143732
$ws.Range("B9").Value = 0.1633529992238181
$ws.Range("D9").Value = 0.01740343617858287
$ws.Range("E9").Value = 0.5368515958397637
$ws.Range("F9").Value = 0.6117683615024561
$ws.Range("G9").Value = 0.002389986962527024
$ws.Range("K9").Value = 0.9230439732951368
$ws.Range("N9").Value = 1.185700871953227
$ws.Range("O9").Value = 1.981105835985005
$ws.Range("B10").Value = 0.1814495247310504
$ws.Range("D10").Value = 0.01950413546148866
$ws.Range("E10").Value = 0.6364051878545922
$ws.Range("F10").Value = 0.6318634149605344
$ws.Range("G10").Value = 0.002385491102005533
$ws.Range("K10").Value = 1.091717562952795
$ws.Range("N10").Value = 1.165018361735406
$ws.Range("O10").Value = 2.022228328876565
$ws.Range("B11").Value = 0.1897510442411345
$ws.Range("D11").Value = 0.0204536996108402
$ws.Range("E11").Value = 0.6818597429769255
$ws.Range("F11").Value = 0.641494400929048
$ws.Range("G11").Value = 0.002383541581280501
$ws.Range("K11").Value = 1.168269118231649
$ws.Range("N11").Value = 1.156136437074714
$ws.Range("O11").Value = 2.04255388906293
$ws.Range("B12").Value = 0.1929044336439603
$ws.Range("D12").Value = 0.02081238365921934
$ws.Range("E12").Value = 0.6990984027261362
$ws.Range("F12").Value = 0.6452122308921702
$ws.Range("G12").Value = 0.002382817026465679
$ws.Range("K12").Value = 1.197231095560937
$ws.Range("N12").Value = 1.152848819699706
$ws.Range("O12").Value = 2.050484890008022
$ws.Range("B13").Value = 0.1922248625933207
$ws.Range("D13").Value = 0.02073517484028287
$ws.Range("E13").Value = 0.6953845606395959
$ws.Range("F13").Value = 0.6444083755434349
$ws.Range("G13").Value = 0.002382972464654292
$ws.Range("K13").Value = 1.190994802464274
$ws.Range("N13").Value = 1.15355349609667
$ws.Range("O13").Value = 2.048766370974505
$ws.Range("B14").Value = 0.190010280450764
$ws.Range("D14").Value = 0.02048322684726145
$ws.Range("E14").Value = 0.6832774453084198
$ws.Range("F14").Value = 0.6417988482583752
$ws.Range("G14").Value = 0.002383481697872911
$ws.Range("K14").Value = 1.170652371371148
$ws.Range("N14").Value = 1.155864443916549
$ws.Range("O14").Value = 2.043201676581276
$ws.Range("B15").Value = 0.1886550536651725
$ws.Range("D15").Value = 0.02032878410609129
$ws.Range("E15").Value = 0.6758649332800104
$ws.Range("F15").Value = 0.6402096668321207
$ws.Range("G15").Value = 0.002383795397979191
$ws.Range("K15").Value = 1.158188581397098
$ws.Range("N15").Value = 1.157289835629996
$ws.Range("O15").Value = 2.039823675272231
$ws.Range("B16").Value = 0.1809083800734896
$ws.Range("D16").Value = 0.01944195526657921
$ws.Range("E16").Value = 0.6334381570510317
$ws.Range("F16").Value = 0.6312438922337549
$ws.Range("G16").Value = 0.002385620426824418
$ws.Range("K16").Value = 1.086711077965788
$ws.Range("N16").Value = 1.165609417553604
$ws.Range("O16").Value = 2.020932696572856
$ws.Range("B17").Value = 0.1761736702885486
$ws.Range("D17").Value = 0.01889634618341773
$ws.Range("E17").Value = 0.6074549071603172
$ws.Range("F17").Value = 0.6258693554453743
$ws.Range("G17").Value = 0.002386764475610685
$ws.Range("K17").Value = 1.04281561579063
$ws.Range("N17").Value = 1.170848129051109
$ws.Range("O17").Value = 2.009759268810939
$ws.Range("B18").Value = 0.173456926115378
$ws.Range("D18").Value = 0.01858195785191441
$ws.Range("E18").Value = 0.5925256406874411
$ws.Range("F18").Value = 0.6228241292676699
$ws.Range("G18").Value = 0.002387431511302085
$ws.Range("K18").Value = 1.017551309293765
$ws.Range("N18").Value = 1.173910865573689
$ws.Range("O18").Value = 2.003484888269497
$ws.Range("B19").Value = 0.1725382113248202
$ws.Range("D19").Value = 0.01847541449363632
$ws.Range("E19").Value = 0.5874734690612655
$ws.Range("F19").Value = 0.6218009683925487
$ws.Range("G19").Value = 0.002387658907383564
$ws.Range("K19").Value = 1.00899438591955
$ws.Range("N19").Value = 1.174956367702464
$ws.Range("O19").Value = 2.001386603080419
$ws.Range("B20").Value = 0.176677013143248
$ws.Range("D20").Value = 0.01895448621018403
$ws.Range("E20").Value = 0.6102192385104246
$ws.Range("F20").Value = 0.6264367137757034
$ws.Range("G20").Value = 0.002386641757823686
$ws.Range("K20").Value = 1.04749010632969
$ws.Range("N20").Value = 1.170285328843413
$ws.Range("O20").Value = 2.010932928916958
$ws.Range("B21").Value = 0.1906604928369546
$ws.Range("D21").Value = 0.02055725457307744
$ws.Range("E21").Value = 0.6868328769331811
$ws.Range("F21").Value = 0.6425634056115257
$ws.Range("G21").Value = 0.002383331752756378
$ws.Range("K21").Value = 1.176628159546283
$ws.Range("N21").Value = 1.155183605657548
$ws.Range("O21").Value = 2.044829796336558
$ws.Range("B22").Value = 0.1998564490841943
$ws.Range("D22").Value = 0.02159952664533904
$ws.Range("E22").Value = 0.7370568330484844
$ws.Range("F22").Value = 0.6535158929295761
$ws.Range("G22").Value = 0.002381248215787668
$ws.Range("K22").Value = 1.26087293419323
$ws.Range("N22").Value = 1.145755452531624
$ws.Range("O22").Value = 2.068348825309016
$ws.Range("B23").Value = 0.1949432446113093
$ws.Range("D23").Value = 0.02104373303431117
$ws.Range("E23").Value = 0.7102367499647499
$ws.Range("F23").Value = 0.6476324524003587
$ws.Range("G23").Value = 0.002382352965526803
$ws.Range("K23").Value = 1.215924272999928
$ws.Range("N23").Value = 1.150747004871974
$ws.Range("O23").Value = 2.055670878149414
$ws.Range("B24").Value = 0.1764494353052442
$ws.Range("D24").Value = 0.01892820331437406
$ws.Range("E24").Value = 0.6089694569288469
$ws.Range("F24").Value = 0.6261800719640291
$ws.Range("G24").Value = 0.002386697209564304
$ws.Range("K24").Value = 1.045376857016095
$ws.Range("N24").Value = 1.170539612164255
$ws.Range("O24").Value = 2.010401851979651
$ws.Range("B25").Value = 0.1567514033238382
$ws.Range("D25").Value = 0.01662462404203069
$ws.Range("E25").Value = 0.5003249013319788
$ws.Range("F25").Value = 0.6047997540993748
$ws.Range("G25").Value = 0.002391727911860369
$ws.Range("K25").Value = 0.8607739112794661
$ws.Range("N25").Value = 1.193784146755213
$ws.Range("O25").Value = 1.967388019771164
